$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("~dicom_tag_dumps")

# Fill in the previously-empty transferSyntaxUid column (I) for rows 2 and 3
$ws.Range("I2").Value = "LittleEndianImplicit"
$ws.Range("I3").Value = "JPEGLossless:Non-hierarchical-1stOrderPrediction"

# Widen column I to fit the new, longer content
# (closest achievable value to the target 60.7109375 stored width, since
# Excel's ColumnWidth property quantizes to whole-pixel steps)
$ws.Columns.Item(9).ColumnWidth = 59.8
